$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Paragraph 1: "This is a Microsoft word document." gains a
#    trailing two spaces, then three red (C00000) runs spelling out
#    "(This is a change \u2013 Version for branch alternate)".
# ------------------------------------------------------------------

# Add the two trailing spaces to the existing (uncoloured) run.
$d.Paragraphs(1).Range.InsertAfter("  ")

# Run 1 (coloured)
$p1 = $d.Paragraphs(1).Range
$startPos = $p1.End - 1
$txt1 = "(This is a change " + [char]0x2013 + " Ve"
$p1.InsertAfter($txt1)
$run1 = $d.Range($startPos, $startPos + $txt1.Length)
$run1.Font.Color = 192

# Run 2 (coloured)
$p1 = $d.Paragraphs(1).Range
$startPos = $p1.End - 1
$txt2 = "rsion for branch alternate"
$p1.InsertAfter($txt2)
$run2 = $d.Range($startPos, $startPos + $txt2.Length)
$run2.Font.Color = 192

# Run 3 (coloured)
$p1 = $d.Paragraphs(1).Range
$startPos = $p1.End - 1
$txt3 = ")"
$p1.InsertAfter($txt3)
$run3 = $d.Range($startPos, $startPos + $txt3.Length)
$run3.Font.Color = 192

# ------------------------------------------------------------------
# 2) Insert a new, empty paragraph right after
#    "It will be treated as a binary file by Git." carrying the
#    shading + paragraph-mark run formatting from the diff.
# ------------------------------------------------------------------

$gitPara = $d.Paragraphs(2).Range
$gitPara.InsertParagraphAfter()

$newPara = $d.Paragraphs(3).Range

$newParaXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:color w:val="202122"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$newPara.InsertXML($newParaXml)

Write-Output "edit complete"
